$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (pushes rows 13.. down to 14..)
$ws.Rows.Item(13).Insert()

# The inserted row 13 picked up column-A formatting in A13 (from the row above);
# clear it since the new row only has B/C content (no A13 in the target layout).
$ws.Range("A13").Clear()

# Copy the normal column B / column C cell formatting onto the new row 13
# cells before setting their values (the Insert() leaves them using the
# row-above's style otherwise).
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null

# Row 10 (Objetivos:) - fix B/C content
$ws.Range("B10").Value = "Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão"
$ws.Range("C10").Value = "Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;Caracterizar as diversas áreas funcionais existentes nas organizações;Despertar o interesse dos alunos para questões de gestão"

# New row 13 - "Docentes responsáveis" name, only B/C (no A)
$ws.Range("B13").Value = "849935 - Humberto Felipe da Silva"
$ws.Range("C13").Value = "849935 - Humberto Felipe da Silva"

# Row 14 (Programa resumido:) - fix B/C content
$ws.Range("B14").Value = "1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão"
$ws.Range("C14").Value = "1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão"

# Row 16 (Programa:) - fix B/C content
$ws.Range("B16").Value = "1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente."
$ws.Range("C16").Value = "1 - A Administração das organizações - definindo a administração2 - O processo administrativo: planejamento, organização, direção, controle3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente."

# Row 19 (Método:) - fix B/C content
$ws.Range("B19").Value = "O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos."
$ws.Range("C19").Value = "O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos."

# Row 20 (Critério:) - fix B/C content
$ws.Range("B20").Value = "Avaliações em diversos formatos realizadas no decorrer do semestre. O peso maior da avaliação será aplicado ao Seminário Final da Disciplina, quando serão realizadas a apresentação oral do trabalho bem como a entrega do trabalho em formato de artigo; essa avaliação representará 70% da média do semestre."
$ws.Range("C20").Value = "Avaliações em diversos formatos realizadas no decorrer do semestre. O peso maior da avaliação será aplicado ao Seminário Final da Disciplina, quando serão realizadas a apresentação oral do trabalho bem como a entrega do trabalho em formato de artigo; essa avaliação representará 70% da média do semestre."

# Row 21 (Norma de recuperação:) - fix B/C content
$ws.Range("B21").Value = "NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota do trabalho de recuperação"
$ws.Range("C21").Value = "NF = (MF + PR)/ 2 , onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota do trabalho de recuperação"

# Row 22 (Bibliografia:) - fix B/C content
$ws.Range("B22").Value = "LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014."
$ws.Range("C22").Value = "LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014.Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014."

"done"
